# "Fix branch instruction format"
#
# Sheet1 ("Register Design" table): add a new "Flags" register row (row 10).
# Sheet2 ("Instruction Design" table): the beq/bne branch instructions used
#   to share a Source-Register-1/Source-Register-2 pair of operands (plus
#   now-unused beqi/bnei immediate-branch rows). Real branches jump to a
#   Label, so: beq now takes a single Label operand, bne keeps the two
#   source-register operands (moving up into beq's old slot), and the
#   beqi/bnei rows are removed. Every section below (List Ops, Jump) shifts
#   up by two rows as a result.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# Sheet1: new "Flags" register entry (row 10)
# ---------------------------------------------------------------------
$ws1.Cells.Item(10, 1).Value = "Flags"
$ws1.Cells.Item(10, 2).Value = 1
$ws1.Cells.Item(10, 3).Value = "`$flg"
$ws1.Cells.Item(10, 4).Value = "Not directly modifiable by programmer"

# ---------------------------------------------------------------------
# Sheet2: fix the branch instruction rows
# ---------------------------------------------------------------------
# Row 18 is "beq" - give it a single Label operand instead of the two
# source-register operands.
$ws2.Cells.Item(18, 2).Value = "Label"
$ws2.Cells.Item(18, 3).ClearContents()

# Row 19 ("beqi") and what becomes row 20 ("bnei" after the first delete)
# are obsolete immediate-branch variants - remove them entirely. This
# leaves the old "bne" row (previously row 20) in place as the new row 19.
$ws2.Rows.Item(19).Delete()
$ws2.Rows.Item(19).Delete()

# The surviving "bne" row is now the second branch instruction (index 1,
# was index 2 when beqi still existed above it).
$ws2.Cells.Item(19, 5).Value = 1

# ---------------------------------------------------------------------
# View state: keep Sheet1 as the active tab/selection, matching the
# original author's final cursor position on each sheet.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("E14").Select()

$ws1.Activate()
$ws1.Range("D10").Select()
